# Updated cryptos list with GitHub Actions.
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row, and
# for the three rows whose ranking changed, also updates the Coin name (B)
# and Link (C): TheGraph <-> InjectiveProtocol swap ranks (rows 40/41),
# and LidoDAOToken is replaced by Mantle (row 51).
#
# Every value is written with a leading apostrophe so Excel keeps it as
# literal text (these are pre-formatted display strings, e.g. "70.958.94"
# or "  +1.80%  ", not real numbers) instead of silently coercing
# numeric-looking text (like "200.10" or "3.00") into a Number and
# dropping significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.958.94"
$ws.Range("E2").Value = "'  +1.80%  "
$ws.Range("D3").Value = "'3.636.09"
$ws.Range("E3").Value = "'  +3.75%  "
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("D5").Value = "'605.32"
$ws.Range("E5").Value = "'  +0.23%  "
$ws.Range("D6").Value = "'200.10"
$ws.Range("E6").Value = "'  +2.34%  "
$ws.Range("E7").Value = "'  +0.27%  "
$ws.Range("E8").Value = "'  +0.07%  "
$ws.Range("E9").Value = "'  +11.20%  "
$ws.Range("E10").Value = "'  -0.46%  "
$ws.Range("D11").Value = "'54.06"
$ws.Range("E12").Value = "'  +1.98%  "
$ws.Range("E13").Value = "'  +0.91%  "
$ws.Range("D14").Value = "'4.211.88"
$ws.Range("E14").Value = "'  +3.70%  "
$ws.Range("D15").Value = "'649.73"
$ws.Range("E15").Value = "'  +9.27%  "
$ws.Range("E16").Value = "'  +1.57%  "
$ws.Range("D17").Value = "'71.038.94"
$ws.Range("E17").Value = "'  +1.72%  "
$ws.Range("D18").Value = "'3.619.43"
$ws.Range("E18").Value = "'  +3.86%  "
$ws.Range("E19").Value = "'  +0.17%  "
$ws.Range("E20").Value = "'  +0.65%  "
$ws.Range("E21").Value = "'  +1.21%  "
$ws.Range("D22").Value = "'18.70"
$ws.Range("E22").Value = "'  +2.90%  "
$ws.Range("D23").Value = "'5.34"
$ws.Range("E23").Value = "'  +0.64%  "
$ws.Range("D24").Value = "'104.52"
$ws.Range("E24").Value = "'  +2.29%  "
$ws.Range("D26").Value = "'3.00"
$ws.Range("E26").Value = "'  -4.97%  "
$ws.Range("D27").Value = "'10.49"
$ws.Range("E27").Value = "'  -3.22%  "
$ws.Range("D28").Value = "'9.76"
$ws.Range("E28").Value = "'  +2.17%  "
$ws.Range("D29").Value = "'34.07"
$ws.Range("E29").Value = "'  +2.31%  "
$ws.Range("D30").Value = "'4.76"
$ws.Range("E30").Value = "'  +9.99%  "
$ws.Range("D31").Value = "'7.22"
$ws.Range("E32").Value = "'  -1.32%  "
$ws.Range("E33").Value = "'  +0.47%  "
$ws.Range("D34").Value = "'63.43"
$ws.Range("E34").Value = "'  +0.52%  "
$ws.Range("D35").Value = "'4.029.60"
$ws.Range("E35").Value = "'  +8.39%  "
$ws.Range("D36").Value = "'0.0₃0881"
$ws.Range("E36").Value = "'  +6.01%  "
$ws.Range("E37").Value = "'  +0.01%  "
$ws.Range("E38").Value = "'  -1.43%  "
$ws.Range("D39").Value = "'507.43"
$ws.Range("E39").Value = "'  +7.24%  "
$ws.Range("B40").Value = "'InjectiveProtocol"
$ws.Range("C40").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'36.82"
$ws.Range("E40").Value = "'  +1.09%  "
$ws.Range("B41").Value = "'TheGraph"
$ws.Range("C41").Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.390"
$ws.Range("E41").Value = "'  -0.44%  "
$ws.Range("E42").Value = "'  -2.56%  "
$ws.Range("E43").Value = "'  +2.38%  "
$ws.Range("D44").Value = "'0.0461"
$ws.Range("E44").Value = "'  +1.52%  "
$ws.Range("D45").Value = "'3.06"
$ws.Range("E45").Value = "'  +8.64%  "
$ws.Range("D46").Value = "'3.50"
$ws.Range("E46").Value = "'  +6.89%  "
$ws.Range("D47").Value = "'0.140"
$ws.Range("E47").Value = "'  +0.58%  "
$ws.Range("E48").Value = "'  +3.15%  "
$ws.Range("E49").Value = "'  -0.23%  "
$ws.Range("E50").Value = "'  +2.17%  "
$ws.Range("B51").Value = "'Mantle"
$ws.Range("C51").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'1.30"
$ws.Range("E51").Value = "'  +1.33%  "
